$d = $word.ActiveDocument

# The "Soutien a la politique de formation" summary table cell has five
# paragraphs ("Nombre de stagiaires...", "Affaires Maritimes...",
# "LPM/ENSM...", "Etrangers...", "Total de jours/stagiaires...") that were
# justified (w:jc val="both"). Remove that explicit justification so the
# paragraphs fall back to the default alignment (this removes the
# <w:jc w:val="both"/> element from each paragraph's pPr).
$justifiedNeedles = @(
    "Nombre de stagiaires",
    "Affaires Maritimes",
    "LPM/ENSM",
    "rangers :",
    "Total de jours/stagiaires"
)

foreach ($needle in $justifiedNeedles) {
    foreach ($p in $d.Paragraphs) {
        if ($p.Range.Text -like ("*" + $needle + "*")) {
            $p.Format.Alignment = 0
        }
    }
}

# The "Total de jours/stagiaires : /" paragraph ended with two manual
# textWrapping line breaks, followed by an entirely empty paragraph
# (just a paragraph mark). Strip the two trailing breaks and remove the
# now-pointless empty paragraph that trailed them.
$totalPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Total de jours/stagiaires*") {
        $totalPara = $p
    }
}

$breakRange = $d.Range($totalPara.Range.End - 3, $totalPara.Range.End - 1)
$breakRange.Delete()

$totalPara2 = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Total de jours/stagiaires*") {
        $totalPara2 = $p
    }
}
$emptyPara = $totalPara2.Next()
$emptyPara.Range.Delete()
